$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 149; this shifts the old rows 149-153 down to 150-154
# while carrying their formatting (e.g. the date style on column D) along.
$ws.Rows.Item(149).Insert()

# New row 148: brand-new weekly price record (replaces the old row-148 data, which
# has effectively moved down into row 149 below).
$ws.Range("A148").Value = 2
$ws.Range("B148").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C148").Value = 'Coquimbo'
$ws.Range("D148").Value = 44461
$ws.Range("E148").Value = 4
$ws.Range("F148").Value = 100112021
$ws.Range("G148").Value = 'Ají'
$ws.Range("H148").Value = 'Inferno'
$ws.Range("I148").Value = 'Primera'
$ws.Range("J148").Value = 200
$ws.Range("K148").Value = 68000
$ws.Range("L148").Value = 70000
$ws.Range("M148").Value = 69000
$ws.Range("N148").Value = '$/caja 25 kilos'
$ws.Range("O148").Value = 'Provincia de Limarí'
$ws.Range("P148").Value = 2760
$ws.Range("Q148").Value = 25
$ws.Range("R148").Value = 'Hortaliza'

# Row 149 is blank after the insert - restore it with what used to be row 148's
# neighbour (old row 149) so the former row 149 content still exists, just shifted.
$ws.Range("A149").Value = 2
$ws.Range("B149").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C149").Value = 'Coquimbo'
$ws.Range("D149").Value = 44357
$ws.Range("E149").Value = 4
$ws.Range("F149").Value = 100112021
$ws.Range("G149").Value = 'Ají'
$ws.Range("H149").Value = 'Americana (o)'
$ws.Range("I149").Value = 'Primera'
$ws.Range("J149").Value = 100
$ws.Range("K149").Value = 38000
$ws.Range("L149").Value = 40000
$ws.Range("M149").Value = 39000
$ws.Range("N149").Value = '$/caja 25 kilos'
$ws.Range("O149").Value = 'Provincia de Limarí'
$ws.Range("P149").Value = 1560
$ws.Range("Q149").Value = 25
$ws.Range("R149").Value = 'Hortaliza'

# Old row 152 (now shifted to row 153) had its date corrected.
$ws.Range("D153").Value = 44203

Write-Output "done"
